# Wed, May 13, 2020 12:04:32 PM
#
# 1) Slide 6 has a table whose style was re-pointed from the deck's
#    custom "Table_0" style to a different (standard) table-style GUID.
# 2) The deck's two theme parts ("Office Theme" and "Integral") were
#    swapped between the Slide Master and the Notes Master - i.e. the
#    slide design that used to be "Integral" becomes the plain
#    "Office Theme" palette (and vice-versa for the notes side).
#    Table styles must be changed with Table.ApplyStyle(), not by
#    assigning the Style property directly.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 6 -------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{BEE1BB0A-56AF-485F-86FB-B88F1306A1A1}")
    }
}

# --- 2. Swap the "Office Theme" / "Integral" colour schemes ---------
function Set-ThemeColor {
    param($ColorScheme, [int]$Index, [byte]$R, [byte]$G, [byte]$B)
    $ColorScheme.Colors($Index).RGB = [int]$R + ([int]$G * 256) + ([int]$B * 65536)
}

# Apply the plain "Office Theme" palette to the presentation's active
# theme (the one backing the Slide Master); the Notes Master's theme
# carries the complementary "Integral" palette.
$tcs = $p.Slides.Item(1).ThemeColorScheme

Set-ThemeColor $tcs 1  0x00 0x00 0x00   # dk1
Set-ThemeColor $tcs 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor $tcs 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor $tcs 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor $tcs 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor $tcs 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor $tcs 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor $tcs 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor $tcs 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor $tcs 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor $tcs 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor $tcs 12 0x95 0x4F 0x72   # folHlink
